# README2_GP-Lesson-Conventions.docx update
# ------------------------------------------
# Applies the edits described by the commit "update roadmap, teaching mat etc."
#
# wdFindWrap constants used below:
#   0 = wdFindStop, 1 = wdFindContinue, 2 = wdFindAsk
# wdReplace constants:
#   0 = wdReplaceNone, 1 = wdReplaceOne, 2 = wdReplaceAll
# wdBreakType:
#   7 = wdPageBreak

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "... Google Slides Presentations, for maximum readability for teachers"
#    -> "... Google Slides Presentations, where we want to ensure maximum
#         readability for teachers"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Google Slides Presentations, for maximum readability for teachers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Google Slides Presentations, where we want to ensure maximum readability for teachers",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Collapse the split "classroom_" / "7" / "-" / "8" runs (and the
#    equivalent "grades 7-8; delete unnecessary..." runs) back into single
#    runs -- text is unchanged, only the run boundaries are merged.
#    Paragraph 23 (1-based) = "classroom_7-8/  presentations for grades 7-8; ..."
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(23).Range
$p.Find.Execute("classroom_7-8", $true, $false, $false, $false, $false, $true, 1, $false, "classroom_7-8", 2) | Out-Null
$p = $d.Paragraphs(23).Range
$p.Find.Execute("grades 7-8; delete unnecessary folders (e.g. if you're not modifying for G9-12)", $true, $false, $false, $false, $false, $true, 1, $false, "grades 7-8; delete unnecessary folders (e.g. if you're not modifying for G9-12)", 2) | Out-Null

# Paragraph 26 (1-based) = "classroom_9-12/  presentations for grades 9-12; ..."
$p = $d.Paragraphs(26).Range
$p.Find.Execute("classroom_9-12", $true, $false, $false, $false, $false, $true, 1, $false, "classroom_9-12", 2) | Out-Null
$p = $d.Paragraphs(26).Range
$p.Find.Execute("grades 9-12; delete unnecessary folders (e.g. if you're not modifying for G5-6)", $true, $false, $false, $false, $false, $true, 1, $false, "grades 9-12; delete unnecessary folders (e.g. if you're not modifying for G5-6)", 2) | Out-Null

# Paragraph 34 (1-based) = "remote_7-8/  presentations for grades 7-8; ..."
$p = $d.Paragraphs(34).Range
$p.Find.Execute("grades 7-8; delete unnecessary folders (e.g. if you're not modifying for G9-12)", $true, $false, $false, $false, $false, $true, 1, $false, "grades 7-8; delete unnecessary folders (e.g. if you're not modifying for G9-12)", 2) | Out-Null

# Paragraph 37 (1-based) = "remote_9-12/  presentations for grades 9-12; ..."
$p = $d.Paragraphs(37).Range
$p.Find.Execute("grades 9-12; delete unnecessary folders (e.g. if you're not modifying for G5-6)", $true, $false, $false, $false, $false, $true, 1, $false, "grades 9-12; delete unnecessary folders (e.g. if you're not modifying for G5-6)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Collapse the split "... presentation versions" / ", multimedia" / ")"
#    runs for the teaching-materials.xlsx bullet back into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "for keeping track of links associated with lesson materials (esp. different presentation versions, multimedia)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "for keeping track of links associated with lesson materials (esp. different presentation versions, multimedia)",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "assemble_lesson.R" -> "compile-lesson.R" (rename the script), and
#    describe the generated output more precisely.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(65).Range
$p.Find.Execute("assemble_", $true, $false, $false, $false, $false, $true, 1, $false, "compile-", 2) | Out-Null

$p = $d.Paragraphs(65).Range
$p.Find.Execute(
    "generates everything necessary to publish",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "generates everything (JSONs, learning charts, etc) necessary to publish",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Insert a manual page break immediately before the "Best practices:"
#    heading so that heading starts on a fresh page.
# ---------------------------------------------------------------------------
$bp = $d.Paragraphs(68)
$breakRange = $d.Range($bp.Range.Start, $bp.Range.Start)
$breakRange.InsertBreak(7)
$newPara = $d.Paragraphs(68)
$newPara.Range.Font.Name = "Montserrat"
$newPara.Range.Font.Size = 16
$newPara.Range.Font.SizeBi = 16
